$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "N/A" value in B2 (subject-id column for the Control group header row)
$ws.Range("B2").Value = ""

# Manually corrected stress-pattern (W/S) cells per Udofot's review
$ws.Range("F2").Value = "S"
$ws.Range("G2").Value = "W"
$ws.Range("H2").Value = "S"

$ws.Range("E3").Value = "S"
$ws.Range("F3").Value = "S"
$ws.Range("G3").Value = "W"
$ws.Range("H3").Value = "S"
$ws.Range("I3").Value = "W"

$ws.Range("E4").Value = "W"
$ws.Range("F4").Value = "S"
$ws.Range("G4").Value = "W"
$ws.Range("H4").Value = "S"

$ws.Range("F5").Value = "S"
$ws.Range("G5").Value = "W"
$ws.Range("H5").Value = "S"
$ws.Range("L5").Value = "S"

$ws.Range("F6").Value = "S"
$ws.Range("G6").Value = "W"
$ws.Range("H6").Value = "S"

$ws.Range("F7").Value = "S"
$ws.Range("G7").Value = "W"
$ws.Range("H7").Value = "S"
$ws.Range("K7").Value = "W"

$ws.Range("F8").Value = "S"
$ws.Range("G8").Value = "W"
$ws.Range("H8").Value = "S"

$ws.Range("H9").Value = "S"
$ws.Range("J9").Value = "W"

$ws.Range("F10").Value = "S"
$ws.Range("I10").Value = "S"
$ws.Range("J10").Value = "W"
$ws.Range("M10").Value = "S"

$ws.Range("D11").Value = "W"
$ws.Range("I11").Value = "S"

# Update the view: scroll so column D is the left-most visible column, and move the
# active selection to I19 (matches the author's saved cursor position).
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("I19").Select()
